$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> Alvearie Team
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row;
# turn it into the new "Jurisdiction" / "United States of America" row.
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was the second (now redundant) "Contact" row - remove it, shifting
# everything below up by one.
$ws1.Rows.Item(11).Delete()

# --- Sheet 2: "Elements" ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 (the root Extension element): Short/Definition were generic
# placeholders - replace with the StructureDefinition's own title/description.
$ws2.Cells.Item(2, 11).Value = "Episode Body System"
$ws2.Cells.Item(2, 12).Value = "Body system related to the episode of care"
